$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four data rows (2-5) get their Date/Calidad/Volumen/Precios/Unidad values
# rotated between rows, while Mercado/Region/Categoria/Origen/Clasificacion stay put.
# New row 2 <- old row 5, new row 3 <- old row 4, new row 4 <- old row 2, new row 5 <- old row 3

$rows = @(2, 3, 4, 5)
$cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")

# Snapshot current values before overwriting anything
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @{}
    foreach ($c in $cols) {
        $snapshot[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Mapping of new row -> source (old) row
$mapping = @{ 2 = 5; 3 = 4; 4 = 2; 5 = 3 }

foreach ($newRow in $rows) {
    $oldRow = $mapping[$newRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $snapshot[$oldRow][$c]
    }
}
